# Updated symbol list: refresh Price (D) and Volume(1h) (E) columns
# for the changed coin rows, preserving their text cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'327.27"
$ws.Range("E2").Value = "'-0.44%"
$ws.Range("D3").Value = "'44.45"
$ws.Range("E3").Value = "'0.39%"
$ws.Range("D4").Value = "'5.171"
$ws.Range("E4").Value = "'-6.18%"
$ws.Range("E5").Value = "'3.46%"
$ws.Range("D6").Value = "'1.942"
$ws.Range("E6").Value = "'-5.49%"
$ws.Range("D7").Value = "'0.9736"
$ws.Range("E7").Value = "'0.17%"
$ws.Range("E8").Value = "'-4.58%"
$ws.Range("D9").Value = "'0.1141"
$ws.Range("E9").Value = "'2.29%"
$ws.Range("D10").Value = "'0.1905"
$ws.Range("E10").Value = "'1.13%"
$ws.Range("D11").Value = "'0.09699"
$ws.Range("E11").Value = "'-2.64%"
$ws.Range("D12").Value = "'0.04617"
$ws.Range("E12").Value = "'-2.00%"
$ws.Range("D13").Value = "'0.1060"
$ws.Range("E13").Value = "'0.45%"
$ws.Range("D14").Value = "'0.001295"
$ws.Range("E14").Value = "'2.77%"
$ws.Range("D15").Value = "'0.005846"
$ws.Range("E15").Value = "'-2.99%"
$ws.Range("D16").Value = "'3.400"
$ws.Range("E16").Value = "'1.80%"
$ws.Range("D17").Value = "'4.445"
$ws.Range("E17").Value = "'0.48%"
$ws.Range("D18").Value = "'0.3362"
$ws.Range("E18").Value = "'1.82%"
$ws.Range("D19").Value = "'8.680"
$ws.Range("E19").Value = "'-14.65%"
$ws.Range("D20").Value = "'0.1362"
$ws.Range("D22").Value = "'0.04148"
$ws.Range("E22").Value = "'1.09%"
$ws.Range("E23").Value = "'-5.42%"
$ws.Range("D24").Value = "'0.004435"
$ws.Range("E24").Value = "'1.05%"
$ws.Range("D25").Value = "'0.0001303"
$ws.Range("E25").Value = "'1.75%"
$ws.Range("E26").Value = "'-20.02%"
$ws.Range("D38").Value = "'0.02749"
$ws.Range("E38").Value = "'2.96%"
$ws.Range("D39").Value = "'0.05658"
$ws.Range("E39").Value = "'0.28%"
$ws.Range("D40").Value = "'0.007862"
$ws.Range("E40").Value = "'3.39%"
$ws.Range("D41").Value = "'0.1414"
$ws.Range("E41").Value = "'0.15%"
$ws.Range("D42").Value = "'0.007310"
$ws.Range("E42").Value = "'-11.19%"
$ws.Range("E43").Value = "'7.89%"
$ws.Range("D44").Value = "'0.007920"
$ws.Range("E44").Value = "'-4.53%"
$ws.Range("D45").Value = "'0.3499"
$ws.Range("D46").Value = "'0.00006853"
$ws.Range("E46").Value = "'-3.29%"
$ws.Range("E47").Value = "'0.40%"
$ws.Range("D48").Value = "'0.003493"
$ws.Range("E48").Value = "'-3.81%"
$ws.Range("E49").Value = "'40.65%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'0.40%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.40%"
